$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the Cell Line header and values first, in row order, so that the
# shared-strings table is populated in the same order as the target file.
$ws.Range("H1").Value = "Cell Line"

$cellLines = @("HK-2","HK-2","HK-2","UMRC6","UMRC6","UMRC6","UOK262","UOK262","UOK262","UOK262","UOK262","UOK262","UOK+DIDS","UOK+DIDS","UOK+DIDS","UOK+DIDS","UOK262","UOK262","UOK262","siRNA_c","siRNA","siRNA")

for ($i = 0; $i -lt $cellLines.Length; $i++) {
    $row = 2 + $i
    $ws.Range("H$row").Value = $cellLines[$i]
}

$ws.Range("I1").Value = "Pyruvate AIC Difference"
$ws.Range("J1").Value = "Intracellular Lactate AIC Difference"
$ws.Range("K1").Value = "Extracellular Lactate Difference"

$ws.Range("H1:K23").Select()
